$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet
$ws.Name = "Sheet1"

# Insert 4 new rows at the top; everything currently on the sheet shifts
# down by 4 rows (rows, merged cells, etc. all move together).
$ws.Rows("1:4").Insert()

# Populate the new header/master-package rows (row 4 is intentionally left blank).
$ws.Range("A1").Value = "MASTER PACKAGE"

$ws.Range("A2").Value = "WesternGlove Centric8 PROD"
$ws.Range("B2").Value = "M12225BVS563:KONRAD"
$ws.Range("C2").Value = "CONSTRUCTION THREAD DETAILS"
$ws.Range("D2").Value = "MASTER"

$ws.Range("A3").Value = "Items"

# Match the formatting (border + wrap + alignment) used throughout the rest
# of the sheet for the newly-written cells by copying it from an existing,
# already-formatted row.
$ws.Range("A5:D5").Copy()
$ws.Range("A1:D3").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Widen columns B and C to fit the new, longer content.
$ws.Columns(2).ColumnWidth = 21.592447916666668
$ws.Columns(3).ColumnWidth = 29.022135416666668
